# [UPD] - Se actualizó el porcentaje de avance del cronograma del proyecto
#
# Updates the "% DE AVANCE" (H column) progress values on the "Cronograma"
# sheet to reflect newly completed / partially-completed tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# Desarrollo del Acta de Constitución -> 100% complete
$ws.Range("H13").Value = 1

# Desarrollo del Plan de Riesgo -> 100% complete
$ws.Range("H14").Value = 1

# Desarrollo del EDT/WBS -> 100% complete
$ws.Range("H15").Value = 1

# Analisis de Negocio -> 50% complete
$ws.Range("H17").Value = 0.5

# Modelo de Situación Actual -> 50% complete
$ws.Range("H18").Value = 0.5

# Keep the selection on the last-touched cell, matching the author's session
$ws.Range("H19").Select()
